$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: append a new sentence/run after "...Has values SHALLOW,
# INTERMEDIATE, and DEEP." describing the new Depth constructor/toColor().
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Has values SHALLOW, INTERMEDIATE, and DEEP.", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find anchor text for change 1"
}
# Collapse to the end of the found text (right before the paragraph mark) and
# insert a brand-new run there.
$r1.Collapse(0)
[void]$r1.InsertAfter("  Contains a constructor to require and define a color on each value.  Contains a toColor() method that returns the color value for the depth.")

# ---------------------------------------------------------------------------
# Change 2: rewrite the EarthquakeMarker paragraph's description of
# colorDetermine() / the quake-title change into the new, shorter text
# describing the removal of colorDetermine().
# ---------------------------------------------------------------------------

# Locate the precise start boundary: right after the (separate) run containing
# just "Depth" and right before " value to set the property...".
$rStart = $d.Content
$foundStart = $rStart.Find.Execute("factory method returns the correct Depth", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if (-not $foundStart) {
    throw "Could not find start anchor text for change 2"
}
$rStart.Collapse(0)
$startPos = $rStart.Start

# Locate the precise end boundary: right after "...depth numeric value." and
# right before "  Finally I added an isRecent...".
$rEnd = $d.Content
$foundEnd = $rEnd.Find.Execute("title to include the depth numeric value.", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not find end anchor text for change 2"
}
$rEnd.Collapse(0)
$oldEnd = $rEnd.Start

# Sanity-check the text we are about to replace.
$oldRange = $d.Range($startPos, $oldEnd)
$expectedOld = " value to set the property.  colorDetermine() is now a switch statement based on that Depth to return the color to use for the marker.  In hindsight, I should have put colorDetermine() into Depth.  I changed the quake title to include the depth numeric value."
if ($oldRange.Text -ne $expectedOld) {
    throw ("Unexpected existing text: [" + $oldRange.Text + "]")
}

$newRun1 = " value to set the property.  "
$newRun2 = "colorDetermine() is completely removed, since each Depth value knows its color the statement pg.fill(depth.toColor()) is all that" + [char]0x2019 + "s needed in drawMarker() now."
$newText = $newRun1 + $newRun2

# Step A: isolate the region [$startPos, $oldEnd] into its own paragraph so
# that later edits can't bleed into the unrelated runs that follow it
# (rightmost split first so $startPos stays valid).
$isolateEnd = $d.Range($oldEnd, $oldEnd)
[void]$isolateEnd.InsertParagraphAfter()
$isolateStart = $d.Range($startPos, $startPos)
[void]$isolateStart.InsertParagraphAfter()

# After both inserts, the isolated paragraph's content lives at
# [$startPos + 1, $oldEnd + 1).
$contentStart = $startPos + 1
$contentEnd = $oldEnd + 1

# Step B: replace the isolated paragraph's text with the new combined text.
$target = $d.Range($contentStart, $contentEnd)
$target.Text = $newText

# Step C: split internally between the two desired runs.
$splitPos = $contentStart + $newRun1.Length
$sr = $d.Range($splitPos, $splitPos)
[void]$sr.InsertParagraphAfter()

# Step D: remove the 3 helper paragraph marks (rightmost first) to merge the
# text back into a single paragraph while keeping the run split created in
# step C (and the original "Depth" run boundary) intact.
$markEnd = $contentStart + $newText.Length + 1
$d.Range($markEnd, $markEnd + 1).Delete()
$d.Range($splitPos, $splitPos + 1).Delete()
$d.Range($startPos, $startPos + 1).Delete()

Write-Output "done"
